$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.874.71'
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("D3").Value = '3.497.41'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.58'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.607'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.492.98'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.30'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.580'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '45.91'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").Value = '4.062.39'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.26'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '609.84'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.503.27'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '69.930.84'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.17'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.870'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.01'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -18.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.49'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.82'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.70'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.55'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.91'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.07'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.86%  '
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '635.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +11.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.87'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.56'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0995'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.67'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0474'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +8.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '56.35'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("D43").Value = '3.317.01'
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.309'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.91'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '32.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").Value = '0.0₃0686'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.129'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.13'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E51").Value = '  -0.01%  '
